$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws1.Range("D8").Value = "Deadband: % top tension"
$ws1.Range("D9").Value = 5
$ws2.Range("B10:C13").ClearContents()
$ws1.Range("D10").Validation.Delete()
